$d = $word.ActiveDocument

$replacements = @(
    @("94÷9=", "43÷2="),
    @("42÷8=", "89÷6="),
    @("70÷3=", "12÷6="),
    @("74÷9=", "41÷2="),
    @("24÷7=", "57÷7="),
    @("69÷9=", "99÷8="),
    @("72÷6=", "36÷2="),
    @("52÷5=", "75÷2="),
    @("90÷6=", "97÷6="),
    @("32÷8=", "56÷4="),
    @("34÷2=", "77÷9="),
    @("81÷3=", "63÷6="),
    @("74÷2=", "34÷9="),
    @("84÷8=", "99÷9="),
    @("60÷7=", "19÷9="),
    @("65÷6=", "83÷9="),
    @("61÷7=", "15÷9="),
    @("36÷9=", "70÷9="),
    @("16÷6=", "68÷9="),
    @("52÷3=", "10÷7="),
    @("91÷7=", "29÷7="),
    @("39÷3=", "10÷2="),
    @("63÷5=", "28÷2="),
    @("28÷4=", "37÷6="),
    @("19÷2=", "24÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
